$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 71.42
$ws.Range("F2").Value = 72.44

$ws.Range("D3").Value = 57.92
$ws.Range("F3").Value = 62.7

$ws.Range("D4").Value = 53.85
$ws.Range("F4").Value = 49.23
